$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 191 ---
# Reuse the existing date-cell style (A190) instead of creating a new one.
$ws.Range("A190").Copy()
$ws.Range("A191").PasteSpecial(-4122)
$ws.Range("A191").Value = 45471.2916666667

$ws.Range("B191").Value = 0
$ws.Range("C191").Value = 5.05000019073486
$ws.Range("D191").Value = 5.05000019073486
$ws.Range("E191").Value = 5.05000019073486
$ws.Range("F191").Value = 5.05000019073486

# Force text storage (shared string) for the adj_close column, matching the
# source data which stores these as strings, then drop the temporary
# "@" number format so no stray style is left behind.
$g191 = $ws.Range("G191")
$g191.NumberFormat = "@"
$g191.Value = "5.05000019073486"
$g191.Style = "Normal"

$ws.Range("H191").Value = "VLC.MI"

# --- Row 192 ---
$ws.Range("A190").Copy()
$ws.Range("A192").PasteSpecial(-4122)
$ws.Range("A192").Value = 45474.6516666667

$ws.Range("B192").Value = 6480
$ws.Range("C192").Value = 5
$ws.Range("D192").Value = 4.59999990463257
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 4.88000011444092

$g192 = $ws.Range("G192")
$g192.NumberFormat = "@"
$g192.Value = "4.88000011444092"
$g192.Style = "Normal"

$ws.Range("H192").Value = "VLC.MI"

$excel.CutCopyMode = 0
